$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (2021-10-04, serial 44473) is inserted as the
# new row 4, pushing the existing rows 4-11 down to rows 5-12.
$ws.Rows(4).Insert()

$ws.Range("A4").Value = 10
$ws.Range("B4").Value = "Vega Modelo de Temuco"
$ws.Range("C4").Value = "La Araucanía"
$ws.Range("D4").Value2 = 44473
$ws.Range("E4").Value = 9
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100108
$ws.Range("H4").Value = "Tropicales y subtropicales"
$ws.Range("I4").Value = 100108001
$ws.Range("J4").Value = "Guayaba"
$ws.Range("K4").Value = "Sin especificar"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 120
$ws.Range("N4").Value = 1200
$ws.Range("O4").Value = 1200
$ws.Range("P4").Value = 1200
$ws.Range("Q4").Value = "$/kilo"
$ws.Range("R4").Value = "Región de Arica y Parinacota"
$ws.Range("S4").Value = 1200
$ws.Range("T4").Value = 1
